$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text for BAT54W-HG3-18 part (row 6, column E)
$ws.Range("E6").Value = "BAT54W-HG3-18 (or BAT 63-02V H6327 )"

# Widen column E to fit the new, longer description text
$ws.Columns("E").ColumnWidth = 35.5

# Update the active selection to F18 (as if user scrolled/selected that cell before save)
[void]$ws.Range("F18").Select()
